# Insert a new weekly price record for Palta (avocado) as row 37,
# pushing the existing rows 37-75 down to 38-76.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 37 (shifts rows 37:75 -> 38:76).
$ws.Rows.Item(37).Insert()

# Populate the newly inserted row 37 with the new record.
$ws.Cells.Item(37, 1).Value  = 1
$ws.Cells.Item(37, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(37, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(37, 4).Value  = 44566
$ws.Cells.Item(37, 5).Value  = 15
$ws.Cells.Item(37, 6).Value  = "Fruta"
$ws.Cells.Item(37, 7).Value  = 100106
$ws.Cells.Item(37, 8).Value  = "Oleaginosos"
$ws.Cells.Item(37, 9).Value  = 100106002
$ws.Cells.Item(37, 10).Value = "Palta"
$ws.Cells.Item(37, 11).Value = "Fuerte"
$ws.Cells.Item(37, 12).Value = "Tercera"
$ws.Cells.Item(37, 13).Value = 200
$ws.Cells.Item(37, 14).Value = 55000
$ws.Cells.Item(37, 15).Value = 57000
$ws.Cells.Item(37, 16).Value = 56000
$ws.Cells.Item(37, 17).Value = "$/caja 25 kilos"
$ws.Cells.Item(37, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(37, 19).Value = 2240
$ws.Cells.Item(37, 20).Value = 25
